$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3312378
$ws.Range("I6").Value = 5714966
$ws.Range("J6").Value = 909790
$ws.Range("K6").Value = 17144898
$ws.Range("L6").Value = 2729370
$ws.Range("M6").Value = -17144786
$ws.Range("N6").Value = -2729594

$ws.Range("H70").Value = 1298.1154
$ws.Range("I70").Value = 1330.0416
$ws.Range("J70").Value = 915
$ws.Range("K70").Value = 3990.1248
$ws.Range("L70").Value = 2745
$ws.Range("M70").Value = -3720.1248
$ws.Range("N70").Value = -3285

$ws.Range("H73").Value = 1298.1154
$ws.Range("I73").Value = 1330.0416
$ws.Range("J73").Value = 915
$ws.Range("K73").Value = 3990.1248
$ws.Range("L73").Value = 2745
$ws.Range("M73").Value = -3054.1248
$ws.Range("N73").Value = -4617

$ws.Range("H100").Value = 2076.9092
$ws.Range("J100").Value = 2399.2
$ws.Range("L100").Value = 2399.2
$ws.Range("N100").Value = -3481.2

$ws.Range("H107").Value = 882.1539
$ws.Range("I107").Value = 887.3333
$ws.Range("J107").Value = 820
$ws.Range("K107").Value = 887.3333
$ws.Range("L107").Value = 820
$ws.Range("M107").Value = 1032.6667
$ws.Range("N107").Value = -4660

$ws.Range("H132").Value = 5168.95
$ws.Range("I132").Value = 4488.4146
$ws.Range("K132").Value = 13465.2438
$ws.Range("M132").Value = -10935.2438

$ws.Range("H137").Value = 996.03125
$ws.Range("I137").Value = 969.35486
$ws.Range("J137").Value = 1021.0909
$ws.Range("K137").Value = 2908.06458
$ws.Range("L137").Value = 3063.2727
$ws.Range("M137").Value = -358.0645800000002
$ws.Range("N137").Value = -8163.2727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7578.15
$ws.Range("I32").Value = 6207.703
$ws.Range("J32").Value = 21434.889
$ws.Range("K32").Value = 6207.703
$ws.Range("L32").Value = 21434.889
$ws.Range("M32").Value = -5920.703
$ws.Range("N32").Value = -22008.889

$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2000
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("M35").Value = -1594

$ws.Range("H45").Value = 999.75
$ws.Range("I45").Value = 999.75
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 999.75
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -622.75
$ws.Range("N45").ClearContents()

$ws.Range("H61").Value = 1543.5532
$ws.Range("I61").Value = 1057.6818
$ws.Range("K61").Value = 1057.6818
$ws.Range("M61").Value = -845.6818000000001

$ws.Range("H95").Value = 33333
$ws.Range("J95").Value = 33333
$ws.Range("L95").Value = 33333
$ws.Range("N95").Value = -38825

$ws.Range("H132").Value = 2020245.9
$ws.Range("I132").Value = 3222.2593
$ws.Range("J132").Value = 3576235.5
$ws.Range("K132").Value = 9666.777900000001
$ws.Range("L132").Value = 10728706.5
$ws.Range("M132").Value = -7136.777900000001
$ws.Range("N132").Value = -10733766.5

$ws.Range("H136").Value = 1543.5532
$ws.Range("I136").Value = 1057.6818
$ws.Range("K136").Value = 3173.0454
$ws.Range("M136").Value = -623.0454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2523.75
$ws.Range("I99").Value = 2250
$ws.Range("J99").Value = 2688
$ws.Range("K99").Value = 2250
$ws.Range("L99").Value = 2688
$ws.Range("M99").Value = -752
$ws.Range("N99").Value = -5684

$ws.Range("H105").Value = 2597.25
$ws.Range("I105").Value = 2400
$ws.Range("J105").Value = 2794.5
$ws.Range("K105").Value = 2400
$ws.Range("L105").Value = 2794.5
$ws.Range("M105").Value = -653
$ws.Range("N105").Value = -6288.5

$ws.Range("H107").Value = 1407.2667
$ws.Range("I107").Value = 1463.2727
$ws.Range("J107").Value = 1253.25
$ws.Range("K107").Value = 1463.2727
$ws.Range("L107").Value = 1253.25
$ws.Range("M107").Value = 456.7273
$ws.Range("N107").Value = -5093.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1716
$ws.Range("I94").Value = 1798
$ws.Range("J94").Value = 1699.6
$ws.Range("K94").Value = 1798
$ws.Range("L94").Value = 1699.6
$ws.Range("M94").Value = -1347
$ws.Range("N94").Value = -2601.6

$ws.Range("H107").Value = 1693.2
$ws.Range("I107").Value = 359.8
$ws.Range("J107").Value = 2359.9
$ws.Range("K107").Value = 359.8
$ws.Range("L107").Value = 2359.9
$ws.Range("M107").Value = 1560.2
$ws.Range("N107").Value = -6199.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1108.4375
$ws.Range("I132").Value = 1023.5
$ws.Range("J132").Value = 1193.375
$ws.Range("K132").Value = 9211.5
$ws.Range("L132").Value = 10740.375
$ws.Range("M132").Value = -6681.5
$ws.Range("N132").Value = -15800.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 31650.4
$ws.Range("J104").Value = 31650.4
$ws.Range("L104").Value = 31650.4
$ws.Range("N104").Value = -38638.4

$ws.Range("H123").Value = 15050.833
$ws.Range("J123").Value = 15050.833
$ws.Range("L123").Value = 15050.833
$ws.Range("N123").Value = -19950.833

$ws.Range("H132").Value = 3045.081
$ws.Range("I132").Value = 2217.5454
$ws.Range("J132").Value = 3395.1924
$ws.Range("K132").Value = 6652.6362
$ws.Range("L132").Value = 10185.5772
$ws.Range("M132").Value = -4122.6362
$ws.Range("N132").Value = -15245.5772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 588756.2
$ws.Range("I46").Value = 492.08334
$ws.Range("J46").Value = 2000590
$ws.Range("K46").Value = 492.08334
$ws.Range("L46").Value = 2000590
$ws.Range("M46").Value = -304.08334
$ws.Range("N46").Value = -2000966

$ws.Range("H68").Value = 1852
$ws.Range("I68").Value = 1707.6923
$ws.Range("J68").Value = 2790
$ws.Range("K68").Value = 1707.6923
$ws.Range("L68").Value = 2790
$ws.Range("M68").Value = -958.6922999999999
$ws.Range("N68").Value = -4288

$ws.Range("H71").Value = 1852
$ws.Range("I71").Value = 1707.6923
$ws.Range("J71").Value = 2790
$ws.Range("K71").Value = 8538.461499999999
$ws.Range("L71").Value = 13950
$ws.Range("M71").Value = -4794.461499999999
$ws.Range("N71").Value = -21438

$ws.Range("H111").Value = 29591.75
$ws.Range("J111").Value = 29591.75
$ws.Range("L111").Value = 29591.75
$ws.Range("N111").Value = -37771.75

$ws.Range("H136").Value = 1678.0444
$ws.Range("I136").Value = 1415.0741
$ws.Range("J136").Value = 2072.5
$ws.Range("K136").Value = 4245.2223
$ws.Range("L136").Value = 6217.5
$ws.Range("M136").Value = -1695.2223
$ws.Range("N136").Value = -11317.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 8478.333000000001
$ws.Range("J104").Value = 8478.333000000001
$ws.Range("L104").Value = 8478.333000000001
$ws.Range("N104").Value = -15466.333

$ws.Range("H136").Value = 1213
$ws.Range("I136").Value = 1036.0769
$ws.Range("J136").Value = 1826.3334
$ws.Range("K136").Value = 3108.2307
$ws.Range("L136").Value = 5479.0002
$ws.Range("M136").Value = -558.2307000000001
$ws.Range("N136").Value = -10579.0002
